$wb = $excel.ActiveWorkbook

# --- Sheet "Metadata" updates ---
$meta = $wb.Worksheets.Item("Metadata")

# URL: https://hl7.fr/fhir/fr/medication/ValueSet/... -> https://hl7.fr/ig/fhir/medication/ValueSet/...
$meta.Range("B2").Value = "https://hl7.fr/ig/fhir/medication/ValueSet/fr-medication-reconciliation-status"

# Name: FrMedicationReconciliationStatus -> FRMedicationReconciliationStatus
$meta.Range("B4").Value = "FRMedicationReconciliationStatus"

# Title: InterOp'Sante -> Interop'Sante
$meta.Range("B5").Value = "value set Interop'Santé - Statut d'une ligne de traitement d'une FCT"

# Date: 2025-04-10T15:35:36+00:00 -> 2026-01-15T08:54:26+00:00
$meta.Range("B8").Value = "2026-01-15T08:54:26+00:00"

# Jurisdiction: (empty) -> FRANCE
$meta.Range("B11").Value = "FRANCE"

# --- Sheet "Include #0" updates ---
$inc = $wb.Worksheets.Item("Include #0")

# System URI value: https://hl7.fr/fhir/fr/medication/CodeSystem/... -> https://hl7.fr/ig/fhir/medication/CodeSystem/...
$inc.Range("B4").Value = "https://hl7.fr/ig/fhir/medication/CodeSystem/fr-medication-reconciliation-status"
